$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the summary title (new reporting period) ---
# A1 already holds a text string; a plain .Value assignment keeps it text
# (no date-looking pattern here), so this is safe as-is.
$ws.Range("A1").Value = "Summary report for 11/17/2019 through 11/30/2019"

# --- Update the start_date / end_date columns (G and H) for rows 5-25 ---
# These cells store plain text (e.g. "11/10/2019") rather than real dates.
# A direct .Value assignment of a date-shaped string gets auto-converted to
# a date serial by the engine (exactly like typing it into a General cell
# in real Excel), which would change the cell's type/style. To keep the
# cells as plain text (matching the original file), stage the text in a
# scratch cell that is explicitly formatted as Text ("@"), copy it, and
# paste-special (values only) into each destination cell - paste-values
# carries just the literal text, not the source's text formatting, so the
# destination cells keep their original (default) style.
$scratch = $ws.Range("Z100")

$scratch.NumberFormat = "@"
$scratch.Value = "11/17/2019"
$scratch.Copy()
for ($r = 5; $r -le 25; $r++) {
    $ws.Cells.Item($r, 7).PasteSpecial(-4163)  # xlPasteValues -> column G (start_date)
}
$scratch.Clear()

$scratch.NumberFormat = "@"
$scratch.Value = "11/30/2019"
$scratch.Copy()
for ($r = 5; $r -le 25; $r++) {
    $ws.Cells.Item($r, 8).PasteSpecial(-4163)  # xlPasteValues -> column H (end_date)
}
$scratch.Clear()

$excel.CutCopyMode = 0

# --- Update the hours values in column C that changed ---
$ws.Range("C5").Value = 38.18
$ws.Range("C6").Value = 64.02
$ws.Range("C7").Value = 33.92
$ws.Range("C8").Value = 18.02
$ws.Range("C11").Value = 8.970000000000001
$ws.Range("C13").Value = 28.03
$ws.Range("C14").Value = 52.18
$ws.Range("C17").Value = 2.48
$ws.Range("C20").Value = 0.33
